$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.915.85"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.287.38"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.03"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.37"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.641"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.11"
$ws.Range("E10").Value = "  -3.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0976"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.11"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.36"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.627.45"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.27"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.872"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.288.37"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.819.48"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.29"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.51"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.68"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.24"
$ws.Range("E24").Value = "  +5.49%  "
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.57"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.70"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.03"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.46"
$ws.Range("E33").Value = "  +5.13%  "
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0820"
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.88"
$ws.Range("E36").Value = "  +7.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.126"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("E38").Value = "  +10.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.77"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.10"
$ws.Range("E41").Value = "  +13.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.34"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.89"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.217"
$ws.Range("E44").Value = "  +7.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.17"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.67"
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.70"
$ws.Range("E49").Value = "  +7.61%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  -1.41%  "
